# Log file updated, with links of Post70
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow the "Table2" ListObject by one row so the table ref / autofilter
# range expand from B10:F79 to B10:F80, matching a new blog-post entry.
$lo = $ws.ListObjects.Item("Table2")
$lo.ListRows.Add() | Out-Null

$row = 80

# Write the new row's values. Order matters for how new strings land in
# the shared-string table: dev.to link, then title, then hashnode link,
# then the two numeric cells.
$ws.Cells.Item($row, 6).Value = "https://dev.to/rahulmishra05/fixed-partition-operating-system-m05-p03-4lmm"
$ws.Cells.Item($row, 3).Value = "Fixed Partition | Operating System - M05 P03"
$ws.Cells.Item($row, 5).Value = "https://programmingport.hashnode.dev/fixed-partition-or-operating-system-m05-p03"
$ws.Cells.Item($row, 2).Value = 70
$ws.Cells.Item($row, 4).Value = 44182

# Match formatting of the row above (number/date format + hyperlink style).
$ws.Range("B79:F79").Copy()
$ws.Range("B80:F80").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Application.GoTo($ws.Range("E80"), $true)
$ws.Range("E80").Select()
